$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with the data that used to live in row 3 (the Oct 26 vs
# Kings XI Punjab match is being dropped; the Oct 7 vs Chennai Super Kings
# match moves up to row 2).
$ws.Range("A2").Value = " Oct 7 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Varun Chakravarthy "

# The numeric-looking columns (runs/balls/4s/6s/sr) are stored as *text* in
# this sheet (see the original "numberStoredAsText" ignoredError). A plain
# .Value assignment of "1" would be auto-coerced to a real number, so use
# the leading-apostrophe force-text trick, then reset the cell style back
# to Normal so no stray quote-prefix formatting is left behind.
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'1"
$ws.Range("I2").Value = "'0"
$ws.Range("J2").Value = "'0"
$ws.Range("K2").Value = "'100.00"
$ws.Range("G2:K2").Style = "Normal"

# Remove the now-duplicate row 3 entirely, shrinking the used range to A1:K2
$ws.Rows("3:3").Delete()
